$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row: species/debtor replace the old species placeholder cols,
#     and the row grows from B1:G1 out to B1:N1 ---
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (index 85) ---
$ws.Range("B2").Value = "現金"
$ws.Range("C2").Value = "高金素梅"
$ws.Range("D2").Value = "陳麗卿新北市泰山區明志路"
$ws.Range("E2").Value = 6000000
$ws.Range("F2").Value = "96年02月06日"
$ws.Range("G2").Value = "借款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
# Quote-prefix forces this to stay literal text instead of being
# auto-parsed into a date serial (it still reads back as plain "2012-04-30").
$ws.Range("J2").Value = "'2012-04-30"
$ws.Range("K2").Value = "高金素梅"
$ws.Range("L2").Value = 926
$ws.Range("M2").Value = "tmpb18e1"
$ws.Range("N2").Value = 85

# --- Row 3 (index 86) ---
$ws.Range("B3").Value = "現金"
$ws.Range("C3").Value = "局金素梅"
$ws.Range("D3").Value = "石旭松新北市泰山區明志路"
$ws.Range("E3").Value = 4000000
$ws.Range("F3").Value = "96年02月06日"
$ws.Range("G3").Value = "借款"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "'2012-04-30"
$ws.Range("K3").Value = "高金素梅"
$ws.Range("L3").Value = 926
$ws.Range("M3").Value = "tmpb18e1"
$ws.Range("N3").Value = 86

# --- Copy formatting onto the newly-populated columns, after the values
#     are in place, so the quote-prefix text trick above doesn't leave a
#     stray number format behind. ---
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)

$excel.CutCopyMode = 0
